$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2520.1428
$ws.Range("I100").Value = 2188.2
$ws.Range("J100").Value = 3350
$ws.Range("K100").Value = 2188.2
$ws.Range("L100").Value = 3350
$ws.Range("M100").Value = -1647.2
$ws.Range("N100").Value = -4432
$ws.Range("H112").Value = 500849.5
$ws.Range("I112").Value = 334466.34
$ws.Range("J112").Value = 999999
$ws.Range("K112").Value = 1003399.02
$ws.Range("L112").Value = 2999997
$ws.Range("M112").Value = -1002291.02
$ws.Range("N112").Value = -3002213
$ws.Range("H137").Value = 2801.6875
$ws.Range("I137").Value = 2404.9167
$ws.Range("K137").Value = 7214.750100000001
$ws.Range("M137").Value = -4664.750100000001
$ws.Range("H138").Value = 6109.7144
$ws.Range("J138").Value = 7048.154
$ws.Range("L138").Value = 21144.462
$ws.Range("N138").Value = -31424.462
$ws.Range("H141").Value = 2491.762
$ws.Range("I141").Value = 2791.4707
$ws.Range("K141").Value = 8374.4121
$ws.Range("M141").Value = -3194.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 90910850
$ws.Range("I61").Value = 90910850
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 90910850
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -90910638
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 30307548
$ws.Range("I74").Value = 31254348
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 31254348
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -31253474
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 30307548
$ws.Range("I77").Value = 31254348
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 156271740
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -156267372
$ws.Range("N77").Value = -58736
$ws.Range("H132").Value = 4002572.2
$ws.Range("I132").Value = 5002471
$ws.Range("J132").Value = 2976.8
$ws.Range("K132").Value = 15007413
$ws.Range("L132").Value = 8930.400000000001
$ws.Range("M132").Value = -15004883
$ws.Range("N132").Value = -13990.4
$ws.Range("H136").Value = 90910850
$ws.Range("I136").Value = 90910850
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 272732550
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -272730000
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28334572
$ws.Range("I134").Value = 36429624
$ws.Range("J134").Value = 1894.5
$ws.Range("K134").Value = 109288872
$ws.Range("L134").Value = 5683.5
$ws.Range("M134").Value = -109286337
$ws.Range("N134").Value = -10753.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3666.1738
$ws.Range("I31").Value = 6115
$ws.Range("J31").Value = 2481.258
$ws.Range("K31").Value = 6115
$ws.Range("L31").Value = 2481.258
$ws.Range("M31").Value = -5820
$ws.Range("N31").Value = -3071.258
$ws.Range("H34").Value = 3666.1738
$ws.Range("I34").Value = 6115
$ws.Range("J34").Value = 2481.258
$ws.Range("K34").Value = 6115
$ws.Range("L34").Value = 2481.258
$ws.Range("M34").Value = -5913
$ws.Range("N34").Value = -2885.258
$ws.Range("H58").Value = 27784520
$ws.Range("I58").Value = 27784520
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 27784520
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -27784317
$ws.Range("N58").ClearContents()
$ws.Range("H99").Value = 8312.723
$ws.Range("I99").Value = 9325.933999999999
$ws.Range("J99").Value = 3246.6667
$ws.Range("K99").Value = 9325.933999999999
$ws.Range("L99").Value = 3246.6667
$ws.Range("M99").Value = -7827.933999999999
$ws.Range("N99").Value = -6242.6667
$ws.Range("H126").Value = 8312.723
$ws.Range("I126").Value = 9325.933999999999
$ws.Range("J126").Value = 3246.6667
$ws.Range("K126").Value = 27977.802
$ws.Range("L126").Value = 9740.000100000001
$ws.Range("M126").Value = -25507.802
$ws.Range("N126").Value = -14680.0001
$ws.Range("H132").Value = 18185100
$ws.Range("I132").Value = 21279768
$ws.Range("J132").Value = 3925.375
$ws.Range("K132").Value = 63839304
$ws.Range("L132").Value = 11776.125
$ws.Range("M132").Value = -63836774
$ws.Range("N132").Value = -16836.125
$ws.Range("H134").Value = 6100293
$ws.Range("I134").Value = 6252550.5
$ws.Range("K134").Value = 18757651.5
$ws.Range("M134").Value = -18755116.5
$ws.Range("H136").Value = 27784520
$ws.Range("I136").Value = 27784520
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 83353560
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -83351010
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7100.721
$ws.Range("J68").Value = 7254.241
$ws.Range("L68").Value = 21762.723
$ws.Range("N68").Value = -23384.723
$ws.Range("H71").Value = 7100.721
$ws.Range("J71").Value = 7254.241
$ws.Range("L71").Value = 65288.169
$ws.Range("N71").Value = -73400.16899999999
$ws.Range("H109").Value = 1440.2222
$ws.Range("I109").Value = 1440.2222
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 4320.6666
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -3280.6666
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3908360.2
$ws.Range("I132").Value = 5002154
$ws.Range("J132").Value = 1954
$ws.Range("K132").Value = 15006462
$ws.Range("L132").Value = 5862
$ws.Range("M132").Value = -15003932
$ws.Range("N132").Value = -10922

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7356169.5
$ws.Range("I132").Value = 8931899
$ws.Range("J132").Value = 2764.6667
$ws.Range("K132").Value = 26795697
$ws.Range("L132").Value = 8294.000100000001
$ws.Range("M132").Value = -26793167
$ws.Range("N132").Value = -13354.0001
$ws.Range("H136").Value = 5074.75
$ws.Range("I136").Value = 5027.5557
$ws.Range("K136").Value = 15082.6671
$ws.Range("M136").Value = -12532.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1794.0667
$ws.Range("I126").Value = 1651.2
$ws.Range("J126").Value = 2079.8
$ws.Range("K126").Value = 4953.6
$ws.Range("L126").Value = 6239.400000000001
$ws.Range("M126").Value = -2483.6
$ws.Range("N126").Value = -11179.4
$ws.Range("H132").Value = 21749572
$ws.Range("I132").Value = 25004308
$ws.Range("J132").Value = 51332.668
$ws.Range("K132").Value = 75012924
$ws.Range("L132").Value = 153998.004
$ws.Range("M132").Value = -75010394
$ws.Range("N132").Value = -159058.004
$ws.Range("H136").Value = 41669690
$ws.Range("I136").Value = 50002524
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 150007572
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -150005022
$ws.Range("N136").Value = -21600
